$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.403.77"
$ws.Range("E2").Value = "  +3.07%  "

$ws.Range("D3").Value = "2.680.59"
$ws.Range("E3").Value = "  +1.57%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "522.66"
$ws.Range("E5").Value = "  +2.04%  "

$ws.Range("E6").Value = "  +2.19%  "

$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").Value = "0.575"
$ws.Range("E8").Value = "  +1.98%  "

$ws.Range("D9").Value = "2.699.13"
$ws.Range("E9").Value = "  +1.15%  "

$ws.Range("D10").Value = "6.46"
$ws.Range("E10").Value = "  +3.31%  "

$ws.Range("E11").Value = "  +0.35%  "

$ws.Range("D12").Value = "0.341"
$ws.Range("E12").Value = "  +2.44%  "

$ws.Range("D14").Value = "3.150.82"
$ws.Range("E14").Value = "  +1.51%  "

$ws.Range("D15").Value = "60.373.85"
$ws.Range("E15").Value = "  +2.98%  "

$ws.Range("D16").Value = "21.25"
$ws.Range("E16").Value = "  +1.99%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0000139"
$ws.Range("E17").Value = "  +1.68%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.704.37"
$ws.Range("E18").Value = "  +1.53%  "

$ws.Range("D19").Value = "351.17"
$ws.Range("E19").Value = "  +2.57%  "

$ws.Range("D20").Value = "4.55"
$ws.Range("E20").Value = "  +0.85%  "

$ws.Range("D21").Value = "10.62"
$ws.Range("E21").Value = "  +2.65%  "

$ws.Range("D22").Value = "6.33"
$ws.Range("E22").Value = "  +4.14%  "

$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  -0.17%  "

$ws.Range("D24").Value = "62.89"
$ws.Range("E24").Value = "  +3.72%  "

$ws.Range("E25").Value = "  +1.49%  "

$ws.Range("D26").Value = "0.168"
$ws.Range("E26").Value = "  +5.48%  "

$ws.Range("D27").Value = "0.994"
$ws.Range("E27").Value = "  +0.10%  "

$ws.Range("D28").Value = "0.0₃0816"
$ws.Range("E28").Value = "  +1.71%  "

$ws.Range("D29").Value = "7.27"
$ws.Range("E29").Value = "  +1.51%  "

$ws.Range("E30").Value = "  +8.05%  "

$ws.Range("E31").Value = "  +0.11%  "

$ws.Range("E32").Value = "  +2.06%  "

$ws.Range("D33").Value = "19.08"
$ws.Range("E33").Value = "  +1.56%  "

$ws.Range("D34").Value = "148.08"
$ws.Range("E34").Value = "  -0.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.30"
$ws.Range("E35").Value = "  +8.06%  "

$ws.Range("D36").Value = "1.25"
$ws.Range("E36").Value = "  +9.77%  "

$ws.Range("D37").Value = "0.949"
$ws.Range("E37").Value = "  -5.25%  "

$ws.Range("D38").Value = "1.55"
$ws.Range("E38").Value = "  +11.32%  "

$ws.Range("D39").Value = "0.876"
$ws.Range("E39").Value = "  +4.09%  "

$ws.Range("D40").Value = "36.91"
$ws.Range("E40").Value = "  +0.95%  "

$ws.Range("E41").Value = "  +1.05%  "

$ws.Range("D42").Value = "282.03"
$ws.Range("E42").Value = "  +1.67%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "19.97"
$ws.Range("E43").Value = "  +2.54%  "

$ws.Range("E44").Value = "  +1.29%  "

$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "0.996"
$ws.Range("E45").Value = "  +0.17%  "

$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "0.611"
$ws.Range("E46").Value = "  -0.81%  "

$ws.Range("D47").Value = "2.135.49"
$ws.Range("E47").Value = "  +7.59%  "

$ws.Range("E48").Value = "  +1.73%  "

$ws.Range("E49").Value = "  +3.98%  "

$ws.Range("E50").Value = "  +2.31%  "

$ws.Range("E51").Value = "  +1.84%  "
